{"js": "// Apply the \"Research Log\" edits:\n//  1. Remove spell/grammar-check proofing marks (w:proofErr) that were\n//     splitting single sentences into multiple runs, by rewriting each\n//     affected paragraph's text as one clean run (content unchanged).\n//  2. Append two new paragraphs (\"October 18, 2018\" + the note about the\n//     test-set change) right before the trailing empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraphs (0-based) whose text is unchanged but whose runs need to be\n// merged back into a single run, clearing away stray proofErr markers.\nconst cleanupIndexes = [3, 10, 11, 13, 14, 15, 18, 19];\n\nfor (const idx of cleanupIndexes) {\n  const para = paragraphs.items[idx];\n  const text = para.text;\n  para.clear();\n  para.insertText(text, Word.InsertLocation.start);\n}\nawait context.sync();\n\n// Insert the two new paragraphs right after the \"...Classification-Random-\n// Data.ipynb\" paragraph (index 19) and before the trailing empty paragraph.\nconst anchor = paragraphs.items[19];\nconst p1 = anchor.insertParagraph(\n  \"October 18, 2018\",\n  Word.InsertLocation.after\n);\np1.insertParagraph(\n  \"Modified the random analysis script so it test on the actual test data set and trains on the whole training set. The test set was run one once and no parameters were changed to adjust to it.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Apply the \"Research Log\" edits:\n#  1. Remove spell/grammar-check proofing marks (proofErr) that were\n#     splitting single sentences into multiple runs, by rewriting each\n#     affected paragraph's text as one clean run (content unchanged).\n#  2. Append two new paragraphs (\"October 18, 2018\" + the note about the\n#     test-set change) right before the trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n# 1-indexed paragraph numbers (COM Paragraphs collection) whose runs were\n# split apart by proofErr spell/grammar markers and need to collapse back\n# into a single clean run, with the text left unchanged.\n$cleanupIndexes = @(4, 11, 12, 14, 15, 16, 19, 20)\n\nforeach ($i in $cleanupIndexes) {\n    $p = $d.Paragraphs.Item($i)\n    $rng = $p.Range\n    $start = $rng.Start\n    $text = $rng.Text\n    # Strip the trailing paragraph-mark character (chr 13) if present so we\n    # don't duplicate it when we reinsert the text.\n    if ($text.Length -gt 0 -and [int][char]$text[$text.Length - 1] -eq 13) {\n        $text = $text.Substring(0, $text.Length - 1)\n    }\n    $rng.Delete()\n    $newRng = $d.Range($start, $start)\n    $newRng.InsertBefore($text + \"`r\")\n}\n\n# 2. Insert the two new paragraphs right after the \"...Classification-\n#    Random-Data.ipynb\" paragraph (item 20) and before the trailing empty\n#    paragraph.\n$anchor = $d.Paragraphs.Item(20)\n$anchor.Range.InsertParagraphAfter()\n\n$p1 = $d.Paragraphs.Item(21)\n$p1.Range.Text = \"October 18, 2018\"\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item(22)\n$p2.Range.Text = \"Modified the random analysis script so it test on the actual test data set and trains on the whole training set. The test set was run one once and no parameters were changed to adjust to it.\"\n"}
